{"js": "// Lattice-multiplication worksheet: refresh every exercise cell with new\n// operands/digits while leaving the \"  ----\" divider line and the cell's\n// run formatting (sz=32) untouched.\n//\n// Each table cell holds a single run with 5 text segments joined by\n// manual line breaks (<w:br/>):\n//   \"AB x CD\" / \"  C    D\" / \"  ----\" / \"A|    |\" / \"B|    |\"\n// We rebuild that exact 5-line payload per cell (vertical-tab joined,\n// which Office.js maps to <w:br/> on insert) and write it back with a\n// single Replace-mode insertText on the cell's whole range, so the\n// existing run properties (<w:rPr><w:sz w:val=\"32\"/></w:rPr>) are kept.\n\nconst newCells = [\n  [\"23 x 68\", \"  6    8\", \"  ----\", \"2|    |\", \"3|    |\"],\n  [\"83 x 71\", \"  7    1\", \"  ----\", \"8|    |\", \"3|    |\"],\n  [\"16 x 30\", \"  3    0\", \"  ----\", \"1|    |\", \"6|    |\"],\n  [\"18 x 40\", \"  4    0\", \"  ----\", \"1|    |\", \"8|    |\"],\n  [\"61 x 62\", \"  6    2\", \"  ----\", \"6|    |\", \"1|    |\"],\n  [\"49 x 77\", \"  7    7\", \"  ----\", \"4|    |\", \"9|    |\"],\n  [\"93 x 53\", \"  5    3\", \"  ----\", \"9|    |\", \"3|    |\"],\n  [\"64 x 61\", \"  6    1\", \"  ----\", \"6|    |\", \"4|    |\"],\n  [\"52 x 30\", \"  3    0\", \"  ----\", \"5|    |\", \"2|    |\"],\n  [\"63 x 96\", \"  9    6\", \"  ----\", \"6|    |\", \"3|    |\"],\n  [\"18 x 77\", \"  7    7\", \"  ----\", \"1|    |\", \"8|    |\"],\n  [\"23 x 16\", \"  1    6\", \"  ----\", \"2|    |\", \"3|    |\"],\n  [\"56 x 62\", \"  6    2\", \"  ----\", \"5|    |\", \"6|    |\"],\n  [\"71 x 71\", \"  7    1\", \"  ----\", \"7|    |\", \"1|    |\"],\n  [\"23 x 91\", \"  9    1\", \"  ----\", \"2|    |\", \"3|    |\"],\n];\n\nconst ROWS = 5;\nconst COLS = 3;\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nlet idx = 0;\nfor (let r = 0; r < ROWS; r++) {\n  for (let c = 0; c < COLS; c++) {\n    const lines = newCells[idx];\n    idx++;\n\n    const cell = table.getCell(r, c);\n    const cellBody = cell.body;\n    const fullRange = cellBody.getRange(\"Whole\");\n    // \\v (vertical tab, U+000B) is how Office.js encodes a manual line\n    // break (<w:br/>) inside insertText payloads.\n    fullRange.insertText(lines.join(\"\\v\"), \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Lattice-multiplication worksheet: refresh every exercise cell with new\n# operands/digits while leaving the \"  ----\" divider line and the cell's\n# run formatting (sz=32) untouched.\n#\n# Each table cell holds a single run with 5 lines joined by manual line\n# breaks (<w:br/>), i.e. Chr(11) / vertical-tab in Range.Text:\n#   \"AB x CD\" / \"  C    D\" / \"  ----\" / \"A|    |\" / \"B|    |\"\n# Setting Cell.Range.Text to a Chr(11)-joined string reproduces that\n# structure and keeps the existing run properties intact.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$vt = [char]11\n\n$newCells = @(\n  @(\"23 x 68\", \"  6    8\", \"  ----\", \"2|    |\", \"3|    |\"),\n  @(\"83 x 71\", \"  7    1\", \"  ----\", \"8|    |\", \"3|    |\"),\n  @(\"16 x 30\", \"  3    0\", \"  ----\", \"1|    |\", \"6|    |\"),\n  @(\"18 x 40\", \"  4    0\", \"  ----\", \"1|    |\", \"8|    |\"),\n  @(\"61 x 62\", \"  6    2\", \"  ----\", \"6|    |\", \"1|    |\"),\n  @(\"49 x 77\", \"  7    7\", \"  ----\", \"4|    |\", \"9|    |\"),\n  @(\"93 x 53\", \"  5    3\", \"  ----\", \"9|    |\", \"3|    |\"),\n  @(\"64 x 61\", \"  6    1\", \"  ----\", \"6|    |\", \"4|    |\"),\n  @(\"52 x 30\", \"  3    0\", \"  ----\", \"5|    |\", \"2|    |\"),\n  @(\"63 x 96\", \"  9    6\", \"  ----\", \"6|    |\", \"3|    |\"),\n  @(\"18 x 77\", \"  7    7\", \"  ----\", \"1|    |\", \"8|    |\"),\n  @(\"23 x 16\", \"  1    6\", \"  ----\", \"2|    |\", \"3|    |\"),\n  @(\"56 x 62\", \"  6    2\", \"  ----\", \"5|    |\", \"6|    |\"),\n  @(\"71 x 71\", \"  7    1\", \"  ----\", \"7|    |\", \"1|    |\"),\n  @(\"23 x 91\", \"  9    1\", \"  ----\", \"2|    |\", \"3|    |\")\n)\n\n$rows = 5\n$cols = 3\n$idx = 0\nfor ($r = 1; $r -le $rows; $r++) {\n  for ($c = 1; $c -le $cols; $c++) {\n    $lines = $newCells[$idx]\n    $idx++\n    $newText = [string]::Join($vt, $lines)\n    $cell = $t.Cell($r, $c)\n    $cell.Range.Text = $newText\n  }\n}\n"}
